$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '95.061.16'
$ws.Range('E2').NumberFormat = "@"
$ws.Range('E2').Value = '  -2.14%  '
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '3.583.06'
$ws.Range('E3').NumberFormat = "@"
$ws.Range('E3').Value = '  -2.97%  '
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '0.999'
$ws.Range('E5').NumberFormat = "@"
$ws.Range('E5').Value = '  +20.57%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '224.69'
$ws.Range('E6').NumberFormat = "@"
$ws.Range('E6').Value = '  -5.25%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '633.54'
$ws.Range('E7').NumberFormat = "@"
$ws.Range('E7').Value = '  -3.79%  '
$ws.Range('E8').NumberFormat = "@"
$ws.Range('E8').Value = '  -3.72%  '
$ws.Range('E9').NumberFormat = "@"
$ws.Range('E9').Value = '  +0.10%  '
$ws.Range('E10').NumberFormat = "@"
$ws.Range('E10').Value = '  +0.07%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '3.580.34'
$ws.Range('E11').NumberFormat = "@"
$ws.Range('E11').Value = '  -3.04%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '46.33'
$ws.Range('E12').NumberFormat = "@"
$ws.Range('E12').Value = '  +5.01%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '0.205'
$ws.Range('E13').NumberFormat = "@"
$ws.Range('E13').Value = '  -1.86%  '
$ws.Range('E14').NumberFormat = "@"
$ws.Range('E14').Value = '  -6.55%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '6.42'
$ws.Range('E15').NumberFormat = "@"
$ws.Range('E15').Value = '  -5.34%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '4.251.60'
$ws.Range('E16').NumberFormat = "@"
$ws.Range('E16').Value = '  -2.91%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '94.875.97'
$ws.Range('E17').NumberFormat = "@"
$ws.Range('E17').Value = '  -2.08%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '8.72'
$ws.Range('E18').NumberFormat = "@"
$ws.Range('E18').Value = '  -4.88%  '
$ws.Range('B19').NumberFormat = "@"
$ws.Range('B19').Value = 'WrappedEther'
$ws.Range('C19').NumberFormat = "@"
$ws.Range('C19').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '3.588.45'
$ws.Range('E19').NumberFormat = "@"
$ws.Range('E19').Value = '  -2.27%  '
$ws.Range('B20').NumberFormat = "@"
$ws.Range('B20').Value = 'Chainlink'
$ws.Range('C20').NumberFormat = "@"
$ws.Range('C20').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '19.64'
$ws.Range('E20').NumberFormat = "@"
$ws.Range('E20').Value = '  +4.87%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '13.03'
$ws.Range('E21').NumberFormat = "@"
$ws.Range('E21').Value = '  +0.33%  '
$ws.Range('E22').NumberFormat = "@"
$ws.Range('E22').Value = '  -0.13%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '497.60'
$ws.Range('E23').NumberFormat = "@"
$ws.Range('E23').Value = '  -4.41%  '
$ws.Range('E24').NumberFormat = "@"
$ws.Range('E24').Value = '  -6.32%  '
$ws.Range('E25').NumberFormat = "@"
$ws.Range('E25').Value = '  +19.01%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '116.94'
$ws.Range('E26').NumberFormat = "@"
$ws.Range('E26').Value = '  +15.27%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '0.0000200'
$ws.Range('E27').NumberFormat = "@"
$ws.Range('E27').Value = '  -4.66%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '6.70'
$ws.Range('E28').NumberFormat = "@"
$ws.Range('E28').Value = '  -3.59%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '3.773.75'
$ws.Range('E29').NumberFormat = "@"
$ws.Range('E29').Value = '  -2.99%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '12.43'
$ws.Range('E30').NumberFormat = "@"
$ws.Range('E30').Value = '  -7.62%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '12.88'
$ws.Range('E31').NumberFormat = "@"
$ws.Range('E31').Value = '  +2.11%  '
$ws.Range('E32').NumberFormat = "@"
$ws.Range('E32').Value = '  -4.77%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '0.999'
$ws.Range('E33').NumberFormat = "@"
$ws.Range('E33').Value = '  +0.00%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '1.00'
$ws.Range('E34').NumberFormat = "@"
$ws.Range('E34').Value = '  -0.08%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '0.177'
$ws.Range('E35').NumberFormat = "@"
$ws.Range('E35').Value = '  -6.60%  '
$ws.Range('E36').NumberFormat = "@"
$ws.Range('E36').Value = '  -7.49%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '31.39'
$ws.Range('E37').NumberFormat = "@"
$ws.Range('E37').Value = '  -2.52%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.579'
$ws.Range('E38').NumberFormat = "@"
$ws.Range('E38').Value = '  -2.11%  '
$ws.Range('E39').NumberFormat = "@"
$ws.Range('E39').Value = '  -0.01%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '586.03'
$ws.Range('E40').NumberFormat = "@"
$ws.Range('E40').Value = '  -10.00%  '
$ws.Range('E41').NumberFormat = "@"
$ws.Range('E41').Value = '  -6.58%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '6.77'
$ws.Range('E42').NumberFormat = "@"
$ws.Range('E42').Value = '  -1.07%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '40.55'
$ws.Range('E43').NumberFormat = "@"
$ws.Range('E43').Value = '  +0.05%  '
$ws.Range('E44').NumberFormat = "@"
$ws.Range('E44').Value = '  -2.78%  '
$ws.Range('E45').NumberFormat = "@"
$ws.Range('E45').Value = '  -7.64%  '
$ws.Range('E46').NumberFormat = "@"
$ws.Range('E46').Value = '  +0.37%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '1.88'
$ws.Range('E47').NumberFormat = "@"
$ws.Range('E47').Value = '  -7.92%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '0.909'
$ws.Range('E48').NumberFormat = "@"
$ws.Range('E48').Value = '  -5.42%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '23.46'
$ws.Range('E49').NumberFormat = "@"
$ws.Range('E49').Value = '  -0.71%  '
$ws.Range('E50').NumberFormat = "@"
$ws.Range('E50').Value = '  +2.81%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '220.40'
$ws.Range('E51').NumberFormat = "@"
$ws.Range('E51').Value = '  +7.24%  '
